$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 8 & 9 (big/small explosion FX entries): mark as Done instead of "room for improvement"
$ws.Range("D8").Value = "Done"
$ws.Range("D9").Value = "Done"

# Row 11: Respawn noise FX now has a filename and is Done
$ws.Range("A11").Value = "321.wav"
$ws.Range("D11").Value = "Done"

# Row 12: in-game theme music now has a filename, reworded description, and is Done
$ws.Range("A12").Value = "Game.wav"
$ws.Range("C12").Value = "theme music (upbeat, fast-paced, 8-bit/synth sounds) -- while in game, ramp up intensity gradually"
$ws.Range("D12").Value = "Done"

# Row 13: game over/high score/main menu music now has a filename, reworded description, and is Done
$ws.Range("A13").Value = "Menus.wav"
$ws.Range("C13").Value = "Game over/high scores/main menu music (slower, 8-bit/synth sounds)"
$ws.Range("D13").Value = "Done"

# Update the saved selection to match where the author last clicked
$ws.Range("E7").Select()
